$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I8").Value = "b"
$ws.Range("J8").Value = "Acknowledge (Backchannel)"
$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"
$ws.Range("I13").Value = "sv"
$ws.Range("J13").Value = "Statement-opinion"
$ws.Range("I33").Value = "b"
$ws.Range("J33").Value = "Acknowledge (Backchannel)"
$ws.Range("I45").Value = "b"
$ws.Range("J45").Value = "Acknowledge (Backchannel)"
$ws.Range("I92").Value = "sv"
$ws.Range("J92").Value = "Statement-opinion"
$ws.Range("I104").Value = "aa"
$ws.Range("J104").Value = "Agree/Accept"
$ws.Range("I113").Value = "ba"
$ws.Range("J113").Value = "Appreciation"
$ws.Range("I120").Value = "sd"
$ws.Range("J120").Value = "Statement-non-opinion"
$ws.Range("I123").Value = "sv"
$ws.Range("J123").Value = "Statement-opinion"
$ws.Range("I133").Value = "sd"
$ws.Range("J133").Value = "Statement-non-opinion"
$ws.Range("I151").Value = "b"
$ws.Range("J151").Value = "Acknowledge (Backchannel)"
$ws.Range("I164").Value = "aa"
$ws.Range("J164").Value = "Agree/Accept"
$ws.Range("I170").Value = "b"
$ws.Range("J170").Value = "Acknowledge (Backchannel)"
$ws.Range("I176").Value = "sd"
$ws.Range("J176").Value = "Statement-non-opinion"
$ws.Range("I184").Value = "b"
$ws.Range("J184").Value = "Acknowledge (Backchannel)"
$ws.Range("I198").Value = "ba"
$ws.Range("J198").Value = "Appreciation"
$ws.Range("I200").Value = "sv"
$ws.Range("J200").Value = "Statement-opinion"
$ws.Range("I204").Value = "b"
$ws.Range("J204").Value = "Acknowledge (Backchannel)"
$ws.Range("I212").Value = "sv"
$ws.Range("J212").Value = "Statement-opinion"
$ws.Range("I237").Value = "sv"
$ws.Range("J237").Value = "Statement-opinion"
$ws.Range("I239").Value = "sd"
$ws.Range("J239").Value = "Statement-non-opinion"
$ws.Range("I240").Value = "aa"
$ws.Range("J240").Value = "Agree/Accept"
$ws.Range("I243").Value = "sd"
$ws.Range("J243").Value = "Statement-non-opinion"
$ws.Range("I246").Value = "%"
$ws.Range("J246").Value = "Uninterpretable"
$ws.Range("I247").Value = "sv"
$ws.Range("J247").Value = "Statement-opinion"
$ws.Range("I258").Value = "sd"
$ws.Range("J258").Value = "Statement-non-opinion"
$ws.Range("I261").Value = "sd"
$ws.Range("J261").Value = "Statement-non-opinion"
$ws.Range("I279").Value = "b"
$ws.Range("J279").Value = "Acknowledge (Backchannel)"
$ws.Range("I297").Value = "b"
$ws.Range("J297").Value = "Acknowledge (Backchannel)"
$ws.Range("I309").Value = "ba"
$ws.Range("J309").Value = "Appreciation"
$ws.Range("I328").Value = "b"
$ws.Range("J328").Value = "Acknowledge (Backchannel)"
$ws.Range("I345").Value = "sv"
$ws.Range("J345").Value = "Statement-opinion"
$ws.Range("I351").Value = "sv"
$ws.Range("J351").Value = "Statement-opinion"
$ws.Range("I354").Value = "sv"
$ws.Range("J354").Value = "Statement-opinion"
$ws.Range("I359").Value = "b"
$ws.Range("J359").Value = "Acknowledge (Backchannel)"
$ws.Range("I360").Value = "sv"
$ws.Range("J360").Value = "Statement-opinion"
$ws.Range("I366").Value = "sv"
$ws.Range("J366").Value = "Statement-opinion"
$ws.Range("I368").Value = "sv"
$ws.Range("J368").Value = "Statement-opinion"
$ws.Range("I382").Value = "ba"
$ws.Range("J382").Value = "Appreciation"
